$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: B1 "Időpont" with same formatting as A1 ("Hiba")
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("B1").Value = "Időpont"

# Data row 2: error message + timestamp, default style
$ws.Range("A2").Value = "403 Client Error: Forbidden for url: https://www.utinform.hu/api/datex2/situation"
$ws.Range("B2").Value = "2026-01-30 17:29:33"
